$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (columns E, F, G hold group-code/group-name/category-name
# before the edit, and need to be rotated so that:
#   new E = old G
#   new F = old E
#   new G = old F
# This applies uniformly to every row, including the header row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $eVal = $ws.Cells.Item($r, 5).Value()
    $fVal = $ws.Cells.Item($r, 6).Value()
    $gVal = $ws.Cells.Item($r, 7).Value()

    $ws.Cells.Item($r, 5).Value = $gVal
    $ws.Cells.Item($r, 6).Value = $eVal
    $ws.Cells.Item($r, 7).Value = $fVal
}
